# Params.xlsx edit — "Add files via upload"
#
# 1. eu_secondary_cost!B3:B28 values all reset to a flat 414000
#    (was a declining series starting at 407790).
# 2. eu_secondary_cost selection moved to D17.
# 3. eu_primary_cost becomes the active/selected sheet (was stocklvl),
#    with its selection moved to C8.
# 4. The workbook scrolls the sheet-tab strip so instalable_capacity
#    (3rd tab) is the first visible tab — best-effort, see note below.

$wb = $excel.ActiveWorkbook

# --- eu_secondary_cost: flatten B3:B28 to 414000, move selection to D17 ---
$wsSecondary = $wb.Worksheets.Item("eu_secondary_cost")
$wsSecondary.Range("B3:B28").Value = 414000
$wsSecondary.Range("D17").Select()

# --- eu_primary_cost: becomes the active sheet, selection -> C8 ---
# (done last / re-asserted so it "wins" as the final active sheet, since
# selecting a range on another sheet can otherwise flip the active sheet)
$wsPrimary = $wb.Worksheets.Item("eu_primary_cost")
$wsPrimary.Activate()
$wsPrimary.Range("C8").Select()
$wsPrimary.Activate()

# --- scroll the workbook's tab strip so the 3rd sheet (instalable_capacity)
#     is first visible tab (bookViews/workbookView@firstSheet="2" in the xml) ---
$excel.ActiveWindow.ScrollWorkbookTabs(2, 1)
